# tracker.xlsx update
#  - Fix the "PefEditNative" typo in the O1 header -> "PdfEditNative"
#    (this merges the duplicate shared-string and renumbers the table,
#     matching every other index shift in the diff automatically).
#  - Mark a batch of newly-completed items ("P"/"L" checkmark columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix header typo: "PefEditNative" -> "PdfEditNative"
$ws.Range("O1").Value = "PdfEditNative"

# Newly completed "PdfEditNative" function rows (column P checkmarks)
$doneP = @("P5","P6","P7","P8","P9","P10","P11","P12","P13","P14","P15","P18","P19","P37")
foreach ($addr in $doneP) {
    $ws.Range($addr).Value = 1
}

# Newly completed items in the K/L pair table
$ws.Range("L31").Value = 1
$ws.Range("L32").Value = 1
